# Apply updated trend values to the "Aggregates" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aggregates")

$ws.Range("E2").Value = 1.3
$ws.Range("F2").Value = -0.2
$ws.Range("I2").Value = 0.8
$ws.Range("L2").Value = -0.7
$ws.Range("N2").Value = 1.6
$ws.Range("P2").Value = 2.169999999999999
$ws.Range("R2").Value = 0.017
$ws.Range("T2").Value = 0.05
$ws.Range("V2").Value = 1.1
$ws.Range("X2").Value = 0.8
$ws.Range("Z2").Value = 0.2
$ws.Range("AB2").Value = 3.7
$ws.Range("E3").Value = 0.2
$ws.Range("F3").Value = -0.1
$ws.Range("I3").Value = -1.6
$ws.Range("L3").Value = 0.1
$ws.Range("N3").Value = 0.9
$ws.Range("O3").Value = 52.2
$ws.Range("P3").Value = 3.549999999999999
$ws.Range("R3").Value = 0.017
$ws.Range("T3").Value = -0.025
$ws.Range("V3").Value = 4.4
$ws.Range("X3").Value = 0.8
$ws.Range("Z3").Value = -0.4
$ws.Range("AA3").Value = 8.48
$ws.Range("AB3").Value = 1.14
$ws.Range("E4").Value = -0.2
$ws.Range("F4").Value = 0.5
$ws.Range("I4").Value = -3.8
$ws.Range("L4").Value = -3.5
$ws.Range("N4").Value = -2.5
$ws.Range("P4").Value = 0.2399999999999999
$ws.Range("R4").Value = -0.134
$ws.Range("T4").Value = 0.04899999999999999
$ws.Range("V4").Value = -0.8
$ws.Range("X4").Value = 0.8
$ws.Range("Z4").Value = 0.6
$ws.Range("AB4").Value = 0.3099999999999999
$ws.Range("E5").Value = 0.2
$ws.Range("I5").Value = -1.5
$ws.Range("K5").Value = 22.20000000000001
$ws.Range("L5").Value = 0.9
$ws.Range("N5").Value = 1.1
$ws.Range("O5").Value = 48.77999999999999
$ws.Range("P5").Value = -2.73
$ws.Range("R5").Value = 0.009999999999999998
$ws.Range("T5").Value = 0.008333333333333333
$ws.Range("V5").Value = -0.1
$ws.Range("X5").Value = 0.4
$ws.Range("Z5").Value = -1.6
$ws.Range("AB5").Value = 1.04
$ws.Range("E6").Value = 0.6
$ws.Range("F6").Value = 0.7
$ws.Range("I6").Value = -6.0
$ws.Range("L6").Value = 1.0
$ws.Range("N6").Value = 0.3
$ws.Range("P6").Value = -1.43
$ws.Range("R6").Value = 0.05
$ws.Range("T6").Value = -0.029
$ws.Range("V6").Value = 1.2
$ws.Range("X6").Value = -0.7
$ws.Range("Z6").Value = 0.7
$ws.Range("AB6").Value = 2.87
$ws.Range("E7").Value = 0.1
$ws.Range("F7").Value = 0.7000000000000001
$ws.Range("I7").Value = 0.5
$ws.Range("L7").Value = 3.0
$ws.Range("N7").Value = 2.6
$ws.Range("P7").Value = -3.000000000000001
$ws.Range("Q7").Value = 0.3179999999999999
$ws.Range("R7").Value = -0.026
$ws.Range("T7").Value = -0.06599999999999999
$ws.Range("V7").Value = -1.4
$ws.Range("X7").Value = -0.6
$ws.Range("Z7").Value = 0.2
$ws.Range("AB7").Value = 0.16
$ws.Range("E8").Value = 0.4
$ws.Range("F8").Value = 0.1
$ws.Range("I8").Value = -4.0
$ws.Range("L8").Value = -3.8
$ws.Range("N8").Value = -0.4
$ws.Range("P8").Value = -0.7599999999999995
$ws.Range("R8").Value = 0.067
$ws.Range("T8").Value = 0.01400000000000001
$ws.Range("V8").Value = 2.5
$ws.Range("X8").Value = 2.0
$ws.Range("Z8").Value = 0.6
$ws.Range("AB8").Value = 3.24
$ws.Range("E9").Value = -0.1
$ws.Range("F9").Value = -0.5
$ws.Range("I9").Value = 1.7
$ws.Range("L9").Value = -3.6
$ws.Range("N9").Value = 2.0
$ws.Range("P9").Value = 6.94
$ws.Range("R9").Value = -0.04666666666666668
$ws.Range("T9").Value = 0.066
$ws.Range("V9").Value = -0.3
$ws.Range("X9").Value = 8.4
$ws.Range("Z9").Value = 5.2
$ws.Range("AB9").Value = -0.5299999999999999
$ws.Range("E10").Value = 0.1
$ws.Range("F10").Value = -0.5
$ws.Range("I10").Value = -1.5
$ws.Range("L10").Value = 0.2
$ws.Range("N10").Value = -1.2
$ws.Range("O10").Value = 51.68000000000001
$ws.Range("P10").Value = 1.849999999999999
$ws.Range("R10").Value = -0.017
$ws.Range("S10").Value = 0.8940000000000001
$ws.Range("T10").Value = 0.086
$ws.Range("V10").Value = 0.6
$ws.Range("X10").Value = -1.4
$ws.Range("Z10").Value = 0.6
$ws.Range("AB10").Value = 1.01
$ws.Range("E11").Value = -0.2
$ws.Range("F11").Value = 0.4999999999999999
$ws.Range("I11").Value = 2.4
$ws.Range("L11").Value = 2.3
$ws.Range("N11").Value = 5.2
$ws.Range("P11").Value = 1.310000000000001
$ws.Range("R11").Value = 0.05
$ws.Range("T11").Value = -0.067
$ws.Range("V11").Value = 2.0
$ws.Range("X11").Value = 0.9
$ws.Range("Z11").Value = 0.6
$ws.Range("AB11").Value = -1.42
$ws.Range("E12").Value = -0.1
$ws.Range("F12").Value = 1.0
$ws.Range("I12").Value = 5.0
$ws.Range("L12").Value = -0.4
$ws.Range("N12").Value = 1.7
$ws.Range("P12").Value = 5.480000000000001
$ws.Range("R12").Value = 0.016
$ws.Range("T12").Value = -0.017
$ws.Range("V12").Value = -1.2
$ws.Range("X12").Value = -0.8
$ws.Range("Z12").Value = 0.2
$ws.Range("AB12").Value = -1.39
$ws.Range("E13").Value = 0.1
$ws.Range("F13").Value = 0.6000000000000001
$ws.Range("I13").Value = -3.2
$ws.Range("L13").Value = -1.0
$ws.Range("N13").Value = 0.9
$ws.Range("O13").Value = 52.61999999999999
$ws.Range("P13").Value = -3.42
$ws.Range("R13").Value = -0.05400000000000001
$ws.Range("T13").Value = -0.183
$ws.Range("V13").Value = -0.4
$ws.Range("X13").Value = 2.4
$ws.Range("Z13").Value = 1.2
$ws.Range("AB13").Value = 3.13
